# GraphBuilder feature: extend family-tree sample data (Sheet1 = people,
# Sheet2 = parent/child relations, Sheet3 = spouse pairs) with the
# "father's side" branch (paternal aunts, their husbands and children),
# and refresh the saved worksheet views accordingly.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet1 ("people": name / birthdate / gender) - append rows 17-26
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Activate()

$people = @(
    @("이봉우", "남"),
    @("큰고모", "여"),
    @("작은고모", "여"),
    @("큰고모부", "남"),
    @("작은고모부", "남"),
    @("김춘식", "남"),
    @("김제식", "남"),
    @("김간식", "남"),
    @("홍성우", "남"),
    @("홍성수", "남")
)

$r = 17
foreach ($p in $people) {
    $ws1.Range("A" + $r).Value = $p[0]
    $ws1.Range("C" + $r).Value = $p[1]
    $r = $r + 1
}

# Column A is best-fit-ish widened (~11 chars) once the sheet gets busier.
$ws1.Columns.Item(1).ColumnWidth = 10.2857142857143

# Freeze the header row and leave the saved selection parked past the data.
$ws1.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws1.Range("E31").Select()

# ---------------------------------------------------------------------
# Sheet2 ("부모"/"자식" parent-child edges) - append rows 21-34
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws2.Activate()

$edges = @(
    @("큰고모", "김춘식"),
    @("큰고모", "김간식"),
    @("큰고모", "김제식"),
    @("작은고모", "홍성우"),
    @("작은고모", "홍성수"),
    @("큰고모부", "김춘식"),
    @("큰고모부", "김간식"),
    @("큰고모부", "김제식"),
    @("작은고모부", "홍성우"),
    @("작은고모부", "홍성수"),
    @("이종수", "큰고모"),
    @("이종수", "작은고모"),
    @("전일분", "큰고모"),
    @("전일분", "작은고모")
)

$r = 21
foreach ($e in $edges) {
    $ws2.Range("A" + $r).Value = $e[0]
    $ws2.Range("B" + $r).Value = $e[1]
    $r = $r + 1
}

$ws2.Columns.Item(1).ColumnWidth = 10.2857142857143

$ws2.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

# ---------------------------------------------------------------------
# Sheet3 ("배우자" spouse pairs) - just move the saved selection.
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Sheet3")
$ws3.Activate()
$ws3.Range("A8").Select()

# ---------------------------------------------------------------------
# Leave Sheet2 as the active tab/window, matching the saved workbook view.
# ---------------------------------------------------------------------
$ws2.Activate()
$ws2.Range("A2").Select()
